$d = $word.ActiveDocument

# Locate the anchor paragraph - the last paragraph of the "Clase 03-06-2021"
# section ("...establecer los nuevos valores.") - and collapse to its end.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "Para mi ejercicio en transferir utilice el método girar() y el setSaldo() para establecer los nuevos valores.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor.Collapse(0)

# 22 blank paragraphs, mirroring the spacing used before the previous
# "Clase ..." heading further up in the document.
for ($i = 0; $i -lt 22; $i++) {
    $anchor.InsertParagraphAfter()
    $anchor.Collapse(0)
}

# New heading paragraph - plain text first, bold/italic/underline applied
# afterwards so the formatting doesn't leak into the following paragraphs.
$anchor.InsertParagraphAfter()
$headingPara = $d.Paragraphs.Last
$headingPara.Range.InsertBefore("Clase 04-06-2021 y 05-06-2021")

# Empty paragraph right after the heading (keeps the heading's formatting).
$afterHeading = $headingPara.Range
$afterHeading.Collapse(0)
$afterHeading.InsertParagraphAfter()
$blankPara = $d.Paragraphs.Last

# Body paragraph 1 (plain formatting).
$afterBlank = $blankPara.Range
$afterBlank.Collapse(0)
$afterBlank.InsertParagraphAfter()
$bodyPara1 = $d.Paragraphs.Last
$bodyPara1.Range.Text = "En estas clases realizamos un proyecto de biblioteca en donde se podían prestar y devolver libros."

# Body paragraph 2 (plain formatting).
$afterBody1 = $bodyPara1.Range
$afterBody1.Collapse(0)
$afterBody1.InsertParagraphAfter()
$bodyPara2 = $d.Paragraphs.Last
$bodyPara2.Range.Text = "Para mi proyecto ocupé una variable más de las pedidas (maxLibros) que me decía la cantidad máxima de libros que tenia cada objeto. Con esto pude controlar la cantidad de libros prestados y devueltos para que no excediera la cantidad existente."

# Apply bold / italic / underline to the heading and the blank paragraph
# that follows it (done last so the body paragraphs stay unformatted).
$headingRange = $headingPara.Range
$headingRange.Font.Bold = 1
$headingRange.Font.BoldBi = 1
$headingRange.Font.Italic = 1
$headingRange.Font.ItalicBi = 1
$headingRange.Font.Underline = 1

$blankRange = $blankPara.Range
$blankRange.Font.Bold = 1
$blankRange.Font.BoldBi = 1
$blankRange.Font.Italic = 1
$blankRange.Font.ItalicBi = 1
$blankRange.Font.Underline = 1
